# Generate Report for Handoff
# Updates the localization-status report: bump status from "In Translation"
# to "Ready for handoff" and refresh the associated generation timestamps,
# then widen the (now longer) status columns to fit the new text.

$wb = $excel.ActiveWorkbook

# --- Overview sheet -------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-23 06:37:14"

# --- zh-cn sheet ------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-23 06:37:10"

# --- de-de sheet ------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-23 06:37:14"

# --- Widen the Status columns so the longer text fits -----------------------
$wsOverview.Columns.Item(5).ColumnWidth = 16.333333333333336
$wsOverview.Columns.Item(6).ColumnWidth = 16.333333333333336
$wsZhCn.Columns.Item(3).ColumnWidth = 16.333333333333336
$wsDeDe.Columns.Item(3).ColumnWidth = 16.333333333333336
